$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to numeric values.
$numericLooking = @("D5", "D6", "D8", "D9", "D12", "D15", "D18", "D19", "D20", "D23", "D24", "D27", "D29", "D30", "D31", "D32", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D50")
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated crypto data (price + 1h volume change), matching the
# refreshed GitHub Actions scrape.
$ws.Range("D2").Value = '58.093.91'
$ws.Range("E2").Value = '  -1.94%  '
$ws.Range("D3").Value = '2.468.73'
$ws.Range("E3").Value = '  -2.30%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '518.76'
$ws.Range("E5").Value = '  -3.54%  '
$ws.Range("D6").Value = '132.23'
$ws.Range("E6").Value = '  -4.30%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '0.557'
$ws.Range("E8").Value = '  -1.81%  '
$ws.Range("D9").Value = '0.0992'
$ws.Range("E9").Value = '  -2.41%  '
$ws.Range("E10").Value = '  -1.22%  '
$ws.Range("E11").Value = '  +0.53%  '
$ws.Range("D12").Value = '0.342'
$ws.Range("E12").Value = '  -1.95%  '
$ws.Range("D13").Value = '2.907.20'
$ws.Range("E13").Value = '  -2.41%  '
$ws.Range("D14").Value = '58.030.63'
$ws.Range("E14").Value = '  -1.85%  '
$ws.Range("D15").Value = '22.07'
$ws.Range("E15").Value = '  -4.72%  '
$ws.Range("E16").Value = '  -2.56%  '
$ws.Range("D17").Value = '2.466.70'
$ws.Range("E17").Value = '  -2.72%  '
$ws.Range("D18").Value = '10.85'
$ws.Range("E18").Value = '  -2.54%  '
$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '4.18'
$ws.Range("E19").Value = '  -2.80%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '319.25'
$ws.Range("E20").Value = '  -1.95%  '
$ws.Range("E21").Value = '  -0.17%  '
$ws.Range("E22").Value = '  -4.28%  '
$ws.Range("D23").Value = '64.16'
$ws.Range("D24").Value = '0.409'
$ws.Range("E24").Value = '  -3.62%  '
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("E26").Value = '  -3.65%  '
$ws.Range("D27").Value = '7.39'
$ws.Range("E27").Value = '  -3.89%  '
$ws.Range("D28").Value = '0.0₃0749'
$ws.Range("E28").Value = '  -3.56%  '
$ws.Range("D29").Value = '6.41'
$ws.Range("E29").Value = '  -5.52%  '
$ws.Range("D30").Value = '1.70'
$ws.Range("E30").Value = '  -4.89%  '
$ws.Range("D31").Value = '165.04'
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").Value = '1.14'
$ws.Range("E32").Value = '  -4.60%  '
$ws.Range("E33").Value = '  -0.10%  '
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").Value = '18.13'
$ws.Range("E35").Value = '  -2.00%  '
$ws.Range("E36").Value = '  -9.19%  '
$ws.Range("D37").Value = '3.99'
$ws.Range("E37").Value = '  -3.41%  '
$ws.Range("D38").Value = '1.49'
$ws.Range("E38").Value = '  -4.05%  '
$ws.Range("D39").Value = '0.795'
$ws.Range("E39").Value = '  -2.93%  '
$ws.Range("D40").Value = '3.49'
$ws.Range("E40").Value = '  -4.52%  '
$ws.Range("D41").Value = '275.52'
$ws.Range("E41").Value = '  -5.10%  '
$ws.Range("D42").Value = '5.00'
$ws.Range("E42").Value = '  -4.49%  '
$ws.Range("D43").Value = '0.592'
$ws.Range("E43").Value = '  -3.11%  '
$ws.Range("D44").Value = '126.66'
$ws.Range("E44").Value = '  -4.55%  '
$ws.Range("D45").Value = '0.0908'
$ws.Range("E45").Value = '  -2.76%  '
$ws.Range("D46").Value = '0.0492'
$ws.Range("E46").Value = '  -3.63%  '
$ws.Range("E47").Value = '  -3.39%  '
$ws.Range("D48").Value = '17.12'
$ws.Range("E48").Value = '  -1.75%  '
$ws.Range("D49").Value = '1.733.62'
$ws.Range("E49").Value = '  -1.61%  '
$ws.Range("D50").Value = '0.972'
$ws.Range("E50").Value = '  -1.64%  '
$ws.Range("E51").Value = '  -2.40%  '
